# Auto-generated from verified diff parse: updates market-data-derived columns (H-N)
# across ALC, ARM, BSM, CRP, CUL, GSM, LTW sheets per scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 12138.667
$ws.Range("I111").Value = 3883.3333
$ws.Range("J111").Value = 20394
$ws.Range("K111").Value = 11649.9999
$ws.Range("L111").Value = 61182
$ws.Range("M111").Value = -8582.999899999999
$ws.Range("N111").Value = -67316

$ws.Range("H113").Value = 3075.5
$ws.Range("J113").Value = 3299.8333
$ws.Range("L113").Value = 3299.8333
$ws.Range("N113").Value = -9807.8333

$ws.Range("H116").Value = 13335797
$ws.Range("I116").Value = 100001000
$ws.Range("J116").Value = 2688.3076
$ws.Range("K116").Value = 100001000
$ws.Range("L116").Value = 2688.3076
$ws.Range("M116").Value = -99997558
$ws.Range("N116").Value = -9572.3076

$ws.Range("H132").Value = 2407.125
$ws.Range("I132").Value = 2132.8823
$ws.Range("J132").Value = 3073.1428
$ws.Range("K132").Value = 6398.646900000001
$ws.Range("L132").Value = 9219.428400000001
$ws.Range("M132").Value = -3868.646900000001
$ws.Range("N132").Value = -14279.4284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 2035
$ws.Range("I5").Value = 2493.75
$ws.Range("J5").Value = 200
$ws.Range("K5").Value = 2493.75
$ws.Range("L5").Value = 200
$ws.Range("M5").Value = -2381.75
$ws.Range("N5").Value = -424

$ws.Range("H74").Value = 1308.6316
$ws.Range("I74").Value = 1419.2307
$ws.Range("J74").Value = 1069
$ws.Range("K74").Value = 1419.2307
$ws.Range("L74").Value = 1069
$ws.Range("M74").Value = -545.2307000000001
$ws.Range("N74").Value = -2817

$ws.Range("H77").Value = 1308.6316
$ws.Range("I77").Value = 1419.2307
$ws.Range("J77").Value = 1069
$ws.Range("K77").Value = 7096.1535
$ws.Range("L77").Value = 5345
$ws.Range("M77").Value = -2728.1535
$ws.Range("N77").Value = -14081

$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 2035
$ws.Range("I4").Value = 2493.75
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 2493.75
$ws.Range("L4").Value = 200
$ws.Range("M4").Value = -2378.75
$ws.Range("N4").Value = -430

$ws.Range("H107").Value = 26331.238
$ws.Range("I107").Value = 31585.766
$ws.Range("J107").Value = 3999.5
$ws.Range("K107").Value = 31585.766
$ws.Range("L107").Value = 3999.5
$ws.Range("M107").Value = -29665.766
$ws.Range("N107").Value = -7839.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2536.8223
$ws.Range("I31").Value = 1314.1515
$ws.Range("J31").Value = 5899.1665
$ws.Range("K31").Value = 1314.1515
$ws.Range("L31").Value = 5899.1665
$ws.Range("M31").Value = -1019.1515
$ws.Range("N31").Value = -6489.1665

$ws.Range("H34").Value = 2536.8223
$ws.Range("I34").Value = 1314.1515
$ws.Range("J34").Value = 5899.1665
$ws.Range("K34").Value = 1314.1515
$ws.Range("L34").Value = 5899.1665
$ws.Range("M34").Value = -1112.1515
$ws.Range("N34").Value = -6303.1665

$ws.Range("H58").Value = 1709.8148
$ws.Range("I58").Value = 1335.5264
$ws.Range("K58").Value = 1335.5264
$ws.Range("M58").Value = -1132.5264

$ws.Range("H105").Value = 2881.6667
$ws.Range("I105").Value = 2881.6667
$ws.Range("K105").Value = 2881.6667
$ws.Range("M105").Value = -1134.6667

$ws.Range("H134").Value = 2277.8215
$ws.Range("I134").Value = 1853.6111
$ws.Range("J134").Value = 3041.4
$ws.Range("K134").Value = 5560.8333
$ws.Range("L134").Value = 9124.200000000001
$ws.Range("M134").Value = -3025.8333
$ws.Range("N134").Value = -14194.2

$ws.Range("H136").Value = 1709.8148
$ws.Range("I136").Value = 1335.5264
$ws.Range("K136").Value = 4006.5792
$ws.Range("M136").Value = -1456.5792

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1995
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 1995
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 17955
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -22855

$ws.Range("H131").Value = 16394526
$ws.Range("I131").Value = 512.1111
$ws.Range("J131").Value = 19231952
$ws.Range("K131").Value = 1536.3333
$ws.Range("L131").Value = 57695856
$ws.Range("M131").Value = 3503.6667
$ws.Range("N131").Value = -57705936

$ws.Range("H132").Value = 1845.619
$ws.Range("I132").Value = 1231.5
$ws.Range("J132").Value = 2091.2666
$ws.Range("K132").Value = 11083.5
$ws.Range("L132").Value = 18821.3994
$ws.Range("M132").Value = -8553.5
$ws.Range("N132").Value = -23881.3994

$ws.Range("H133").Value = 5270.643
$ws.Range("J133").Value = 7215.3335
$ws.Range("L133").Value = 21646.0005
$ws.Range("N133").Value = -31766.0005

$ws.Range("H134").Value = 3259.2917
$ws.Range("I134").Value = 2014.375
$ws.Range("J134").Value = 5749.125
$ws.Range("K134").Value = 6043.125
$ws.Range("L134").Value = 17247.375
$ws.Range("M134").Value = -973.125
$ws.Range("N134").Value = -27387.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H113").Value = 1546.5555
$ws.Range("I113").Value = 1122.5333
$ws.Range("K113").Value = 1122.5333
$ws.Range("M113").Value = 1047.4667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2439.5908
$ws.Range("I68").Value = 1954.25
$ws.Range("J68").Value = 3733.8333
$ws.Range("K68").Value = 1954.25
$ws.Range("L68").Value = 3733.8333
$ws.Range("M68").Value = -1205.25
$ws.Range("N68").Value = -5231.8333

$ws.Range("H71").Value = 2439.5908
$ws.Range("I71").Value = 1954.25
$ws.Range("J71").Value = 3733.8333
$ws.Range("K71").Value = 9771.25
$ws.Range("L71").Value = 18669.1665
$ws.Range("M71").Value = -6027.25
$ws.Range("N71").Value = -26157.1665

$ws.Range("H136").Value = 2335.9714
$ws.Range("I136").Value = 1641.2609
$ws.Range("J136").Value = 3667.5
$ws.Range("K136").Value = 4923.7827
$ws.Range("L136").Value = 11002.5
$ws.Range("M136").Value = -2373.7827
$ws.Range("N136").Value = -16102.5
